$d = $word.ActiveDocument

# --- Locate the paragraph that begins "package xx::yy::zz" ---
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "package*xx::yy::zz*") {
        $targetIndex = $i
        break
    }
}

# --- 1. Insert a brand-new empty paragraph right before it ("use lib /path;") ---
$pkgRange = $d.Paragraphs($targetIndex).Range
$pkgRange.InsertParagraphBefore()

$newParaIndex = $targetIndex
$newPara = $d.Paragraphs($newParaIndex).Range

# Type "u" first ...
$newPara.Text = "u"

# ... then continue typing the rest right after it. Because the edit point
# ("u" | "se lib /path;") is where the cursor was last, Word drops its
# "_GoBack" last-edit bookmark exactly there, splitting the run in two.
$afterU = $d.Range($d.Paragraphs($newParaIndex).Range.End - 1, $d.Paragraphs($newParaIndex).Range.End - 1)
$afterU.InsertAfter("se lib /path;")

$bmPoint = $d.Range($afterU.Start, $afterU.Start)
$d.Bookmarks.Add("_GoBack", $bmPoint)

# --- 2. Split the tail of the "package ..." paragraph's last run into three
#        runs by inserting "|LIB" between "@INC" and "/xx/yy/zz.pm" ---
$pkgParaIndex = $targetIndex + 1
$pkgParaRange = $d.Paragraphs($pkgParaIndex).Range

$incEnd = -1
$searchRange = $pkgParaRange.Duplicate
$found = $searchRange.Find.Execute("@INC/xx/yy/zz.pm", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $incEnd = $searchRange.Start + 4
}

$splitPoint = $d.Range($incEnd, $incEnd)
$splitPoint.InsertAfter("|LIB")

# Force the inserted "|LIB" segment (and the text after it) into their own
# runs, then restore formatting so every run ends up with identical rPr
# (sz/szCs 24) just like the surrounding text.
$libRange = $d.Range($incEnd, $incEnd + 4)
$libRange.Font.Bold = 1
$libRange.Font.Bold = 0

Write-Host "done"
